$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new journal entry for "Day 9" directly below the last used row (row 9)
$ws.Range("A10").Value = 45920
$ws.Range("B10").Value = "Citizen Complaint Response Automation final"
$ws.Range("C10").Value = "Citizen Complaint Response Automation"
$ws.Range("D10").Value = "Citizen Complaint Response Automation final.json"

# Match the date number formatting used by the rest of the date column
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat
$ws.Range("A9").NumberFormat = $ws.Range("A3").NumberFormat

# Move the active selection to D11, mirroring what Excel does after entering
# data in the last row of the table
$ws.Range("D11").Select()
